# Append the latest profit-run row (2025-10-08) to the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date to be stored as plain text (matching the other "MM/DD/YYYY"
# rows in the sheet, which are literal strings rather than real date values),
# then clear the temporary "Text" number format so the new cell ends up
# unstyled - just like its neighbours.
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "10/08/2025"
$ws.Range("A52").ClearFormats()

$ws.Range("B52").Value = 14733.41
